$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace subject columns B:E (rows 1-3) with updated data for
# subjects 15 & 16 (meanEMG legmaxROM update)

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 16.446145798608285
$ws.Range("C2").Value = 11.310159249463084
$ws.Range("D2").Value = 13.240537416765733
$ws.Range("E2").Value = 0.041881807647570213

$ws.Range("B3").Value = 33.936555197983608
$ws.Range("C3").Value = 4.0596357921674269
$ws.Range("D3").Value = 2.1218849380856852
$ws.Range("E3").Value = 1.231805422257537

$ws.Range("B1:E3").Select()
